$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    5.7,
    6.3,
    6.8,
    6,
    5.8,
    6.4,
    5.4,
    5.9,
    6.4,
    4.6,
    6.4,
    6.4,
    5.7,
    6.5,
    5.9,
    6.8,
    7.5,
    6.4,
    5.8,
    4.9,
    7.6,
    6.6,
    8.199999999999999,
    8.1,
    8.4,
    7.9,
    8.199999999999999,
    6.8,
    6.6,
    7.9,
    6.5,
    6.4,
    5.3,
    4.5,
    4.7,
    5.2,
    6.4,
    4,
    5,
    5.9
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
